$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing A10 timestamp (tiny float precision update from the WSL re-run)
$ws.Range("A10").Value = 45875.45855712963

# Append the new row 11 with the latest weather-station reading
$ws.Range("A11").Value = 45875.50018407633
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 23
$ws.Range("D11").Value = 20.44
$ws.Range("E11").Value = 75.47
$ws.Range("F11").Value = 613.48
$ws.Range("G11").Value = 11.61
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "12:00:15"
